$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row.
# Update the value from 45177 (2023-09-08) to 45178 (2023-09-09) for all
# data rows (rows 2 through 141).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
